# Fixed PDQ link click assertions & double-beacon issue
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) RightNav (sheet3): the "PageDetail" column is gone, and the
#    "SectionName" column has been repurposed into a "LinkText" column
#    describing the patient/HP toggle link text.
# ---------------------------------------------------------------------
$rightNav = $wb.Worksheets.Item("RightNav")

$rightNav.Range("C1:C5").EntireColumn.Delete()

$rightNav.Range("B1").Value = "LinkText"
$rightNav.Range("B2").Value = "Description of the Evidence"
$rightNav.Range("B3").Value = "View All Sections"
$rightNav.Range("B4").Value = "Descripción de las pruebas"
$rightNav.Range("B5").Value = "Ver todas las secciones"

$rightNav.Columns.Item(1).ColumnWidth = 46.7109375
$rightNav.Columns.Item(2).ColumnWidth = 25.85546875

# ---------------------------------------------------------------------
# 2) PDQDrugPage (sheet2): scrolled/selected a different cell.
# ---------------------------------------------------------------------
$drugPage = $wb.Worksheets.Item("PDQDrugPage")
$drugPage.Application.ActiveWindow.ScrollRow = 22
$drugPage.Range("A5").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) New sheet: PatientHPToggle - holds the patient/HP toggle link data
#    used by the new double-beacon regression test.
# ---------------------------------------------------------------------
$toggle = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$toggle.Name = "PatientHPToggle"

$toggle.Range("A1").Value = "Path"
$toggle.Range("B1").Value = "LinkText"

$toggle.Range("A2").Value = "/espanol/cancer/causas-prevencion/aspectos-generales-prevencion-paciente-pdq"
$toggle.Range("B2").Value = "vaya a la versión para profesionales de salud"

$toggle.Range("A3").Value = "/espanol/cancer/causas-prevencion/aspectos-generales-prevencion-pro-pdq"
$toggle.Range("B3").Value = "vaya a la versión para pacientes"

$toggle.Range("A4").Value = "/about-cancer/causes-prevention/hp-prevention-overview-pdq"
$toggle.Range("B4").Value = "go to patient version"

$toggle.Range("A5").Value = "/about-cancer/causes-prevention/patient-prevention-overview-pdq"
$toggle.Range("B5").Value = "go to health professional version"

$rightNav.Range("A1:B1").Copy() | Out-Null
$toggle.Range("A1:B1").PasteSpecial(-4122) | Out-Null

$toggle.Columns.Item(1).ColumnWidth = 73.4453125
$toggle.Columns.Item(2).ColumnWidth = 38.3046875

$toggle.Range("A6").Select() | Out-Null

# Restore the originally-active sheet/tab.
$rightNav.Activate() | Out-Null
